$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row - correct marks value
$ws.Range("B11").Value = 5

# Update "Total" row - total correct marks and fraction display
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
